$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "29.196.97"
$cell.ClearFormats()
$cell = $ws.Range("E2")
$cell.NumberFormat = "@"
$cell.Value = "  -0.05%  "
$cell.ClearFormats()

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.825.18"
$cell.ClearFormats()
$cell = $ws.Range("E3")
$cell.NumberFormat = "@"
$cell.Value = "  -0.03%  "
$cell.ClearFormats()

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.004"
$cell.ClearFormats()
$cell = $ws.Range("E4")
$cell.NumberFormat = "@"
$cell.Value = "  +0.29%  "
$cell.ClearFormats()

$cell = $ws.Range("E5")
$cell.NumberFormat = "@"
$cell.Value = "  -0.51%  "
$cell.ClearFormats()

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.5990"
$cell.ClearFormats()
$cell = $ws.Range("E6")
$cell.NumberFormat = "@"
$cell.Value = "  -0.65%  "
$cell.ClearFormats()

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "1.004"
$cell.ClearFormats()
$cell = $ws.Range("E7")
$cell.NumberFormat = "@"
$cell.Value = "  +0.18%  "
$cell.ClearFormats()

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.06943"
$cell.ClearFormats()
$cell = $ws.Range("E8")
$cell.NumberFormat = "@"
$cell.Value = "  -2.60%  "
$cell.ClearFormats()

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.2761"
$cell.ClearFormats()
$cell = $ws.Range("E9")
$cell.NumberFormat = "@"
$cell.Value = "  -1.91%  "
$cell.ClearFormats()

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "23.44"
$cell.ClearFormats()
$cell = $ws.Range("E10")
$cell.NumberFormat = "@"
$cell.Value = "  -2.70%  "
$cell.ClearFormats()

$cell = $ws.Range("E11")
$cell.NumberFormat = "@"
$cell.Value = "  -0.61%  "
$cell.ClearFormats()

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "1.876.87"
$cell.ClearFormats()
$cell = $ws.Range("E12")
$cell.NumberFormat = "@"
$cell.Value = "  +1.51%  "
$cell.ClearFormats()

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "4.731"
$cell.ClearFormats()
$cell = $ws.Range("E13")
$cell.NumberFormat = "@"
$cell.Value = "  -1.18%  "
$cell.ClearFormats()

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "0.6285"
$cell.ClearFormats()
$cell = $ws.Range("E14")
$cell.NumberFormat = "@"
$cell.Value = "  -2.19%  "
$cell.ClearFormats()

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.000009816"
$cell.ClearFormats()
$cell = $ws.Range("E15")
$cell.NumberFormat = "@"
$cell.Value = "  +0.32%  "
$cell.ClearFormats()

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "77.43"
$cell.ClearFormats()
$cell = $ws.Range("E16")
$cell.NumberFormat = "@"
$cell.Value = "  -2.50%  "
$cell.ClearFormats()

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "29.005.85"
$cell.ClearFormats()
$cell = $ws.Range("E17")
$cell.NumberFormat = "@"
$cell.Value = "  -0.47%  "
$cell.ClearFormats()

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "5.539"
$cell.ClearFormats()
$cell = $ws.Range("E18")
$cell.NumberFormat = "@"
$cell.Value = "  -7.98%  "
$cell.ClearFormats()

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "215.77"
$cell.ClearFormats()
$cell = $ws.Range("E19")
$cell.NumberFormat = "@"
$cell.Value = "  -6.52%  "
$cell.ClearFormats()

$cell = $ws.Range("E20")
$cell.NumberFormat = "@"
$cell.Value = "  +0.22%  "
$cell.ClearFormats()

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "11.54"
$cell.ClearFormats()
$cell = $ws.Range("E21")
$cell.NumberFormat = "@"
$cell.Value = "  -2.09%  "
$cell.ClearFormats()

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "6.852"
$cell.ClearFormats()
$cell = $ws.Range("E22")
$cell.NumberFormat = "@"
$cell.Value = "  -2.58%  "
$cell.ClearFormats()

$cell = $ws.Range("E23")
$cell.NumberFormat = "@"
$cell.Value = "  -0.10%  "
$cell.ClearFormats()

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "155.72"
$cell.ClearFormats()
$cell = $ws.Range("E24")
$cell.NumberFormat = "@"
$cell.Value = "  -0.18%  "
$cell.ClearFormats()

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "7.943"
$cell.ClearFormats()
$cell = $ws.Range("E25")
$cell.NumberFormat = "@"
$cell.Value = "  -1.82%  "
$cell.ClearFormats()

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "0.1286"
$cell.ClearFormats()
$cell = $ws.Range("E26")
$cell.NumberFormat = "@"
$cell.Value = "  +0.32%  "
$cell.ClearFormats()

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "16.49"
$cell.ClearFormats()
$cell = $ws.Range("E27")
$cell.NumberFormat = "@"
$cell.Value = "  -1.35%  "
$cell.ClearFormats()

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "0.06415"
$cell.ClearFormats()
$cell = $ws.Range("E28")
$cell.NumberFormat = "@"
$cell.Value = "  -5.81%  "
$cell.ClearFormats()

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "1.417"
$cell.ClearFormats()
$cell = $ws.Range("E29")
$cell.NumberFormat = "@"
$cell.Value = "  -3.32%  "
$cell.ClearFormats()

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "1.441"
$cell.ClearFormats()
$cell = $ws.Range("E30")
$cell.NumberFormat = "@"
$cell.Value = "  -1.23%  "
$cell.ClearFormats()

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "3.818"
$cell.ClearFormats()
$cell = $ws.Range("E31")
$cell.NumberFormat = "@"
$cell.Value = "  +1.43%  "
$cell.ClearFormats()

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "3.779"
$cell.ClearFormats()
$cell = $ws.Range("E32")
$cell.NumberFormat = "@"
$cell.Value = "  -1.22%  "
$cell.ClearFormats()

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "1.093"
$cell.ClearFormats()
$cell = $ws.Range("E33")
$cell.NumberFormat = "@"
$cell.Value = "  -3.48%  "
$cell.ClearFormats()

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "1.722"
$cell.ClearFormats()
$cell = $ws.Range("E34")
$cell.NumberFormat = "@"
$cell.Value = "  +0.44%  "
$cell.ClearFormats()

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "0.6455"
$cell.ClearFormats()
$cell = $ws.Range("E35")
$cell.NumberFormat = "@"
$cell.Value = "  -2.06%  "
$cell.ClearFormats()

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "2.543"
$cell.ClearFormats()
$cell = $ws.Range("E36")
$cell.NumberFormat = "@"
$cell.Value = "  +0.22%  "
$cell.ClearFormats()

$cell = $ws.Range("E37")
$cell.NumberFormat = "@"
$cell.Value = "  -0.50%  "
$cell.ClearFormats()

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.01753"
$cell.ClearFormats()
$cell = $ws.Range("E38")
$cell.NumberFormat = "@"
$cell.Value = "  -0.75%  "
$cell.ClearFormats()

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "6.590"
$cell.ClearFormats()
$cell = $ws.Range("E39")
$cell.NumberFormat = "@"
$cell.Value = "  +0.76%  "
$cell.ClearFormats()

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "1.133.50"
$cell.ClearFormats()
$cell = $ws.Range("E40")
$cell.NumberFormat = "@"
$cell.Value = "  -7.45%  "
$cell.ClearFormats()

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.8922"
$cell.ClearFormats()
$cell = $ws.Range("E41")
$cell.NumberFormat = "@"
$cell.Value = "  -3.62%  "
$cell.ClearFormats()

$cell = $ws.Range("E42")
$cell.NumberFormat = "@"
$cell.Value = "  +0.16%  "
$cell.ClearFormats()

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "1.991.33"
$cell.ClearFormats()
$cell = $ws.Range("E43")
$cell.NumberFormat = "@"
$cell.Value = "  +1.50%  "
$cell.ClearFormats()

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "100.89"
$cell.ClearFormats()
$cell = $ws.Range("E44")
$cell.NumberFormat = "@"
$cell.Value = "  +1.19%  "
$cell.ClearFormats()

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "62.10"
$cell.ClearFormats()
$cell = $ws.Range("E45")
$cell.NumberFormat = "@"
$cell.Value = "  -1.89%  "
$cell.ClearFormats()

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.00000000113"
$cell.ClearFormats()
$cell = $ws.Range("E46")
$cell.NumberFormat = "@"
$cell.Value = "  -2.15%  "
$cell.ClearFormats()

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "1.613"
$cell.ClearFormats()
$cell = $ws.Range("E47")
$cell.NumberFormat = "@"
$cell.Value = "  -1.14%  "
$cell.ClearFormats()

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "8.428"
$cell.ClearFormats()
$cell = $ws.Range("E48")
$cell.NumberFormat = "@"
$cell.Value = "  -0.68%  "
$cell.ClearFormats()

$cell = $ws.Range("E49")
$cell.NumberFormat = "@"
$cell.Value = "  -1.40%  "
$cell.ClearFormats()

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.4537"
$cell.ClearFormats()
$cell = $ws.Range("E50")
$cell.NumberFormat = "@"
$cell.Value = "  -0.52%  "
$cell.ClearFormats()

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "6.363"
$cell.ClearFormats()
$cell = $ws.Range("E51")
$cell.NumberFormat = "@"
$cell.Value = "  -3.38%  "
$cell.ClearFormats()
